# Apply changes described in the diff:
# - Sheet1 "orders": remove data rows 2 and 3, keep only the header row.
# - Sheet2 "customer_info": remove data row 2, add "Shipping Details" header in F1.

$wb = $excel.ActiveWorkbook

# --- Sheet1 ("orders"): delete rows 2 and 3 ---
$ws1 = $wb.Worksheets.Item("orders")
$ws1.Rows.Item(2).Resize(2).Delete()

# --- Sheet2 ("customer_info"): add "Shipping Details" header, delete row 2 ---
$ws2 = $wb.Worksheets.Item("customer_info")
$ws2.Range("F1").Value = "Shipping Details"
$ws2.Columns.Item(6).ColumnWidth = 7.43
$ws2.Rows.Item(2).Delete()
